# "Generate Report for Handoff"
#
# The localization-status report is regenerated: the second tracked file
# (ba5b2a27-c42e-4541-af27-6057b6515bf4.*) has dropped out of the report,
# and the status / timestamps for the remaining tracked file
# (8fc7b040-dee2-421b-8e4f-1d316658501e.*) have moved forward to reflect a
# fresh "Ready for handoff" run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"
$ws1.Range("D2").Value = "2016-36-19 12:36:45"

# Row 3 (the ba5b2a27... file) no longer belongs in the report.
$ws1.Rows.Item(3).Delete()

# Rebuild the hyperlink collection so only the row-2 link survives.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/276dd498861d647be956936a380962ede6cdd123/e2e/8fc7b040-dee2-421b-8e4f-1d316658501e.md", "", "", "8fc7b040-dee2-421b-8e4f-1d316658501e.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("E2").Value = "2016-03-19 12:36:42"

# Row 3 (the ba5b2a27... file) no longer belongs in the report.
$ws2.Rows.Item(3).Delete()

# Rebuild the hyperlink collection so only the row-2 links survive.
$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/276dd498861d647be956936a380962ede6cdd123/e2e/8fc7b040-dee2-421b-8e4f-1d316658501e.md", "", "", "8fc7b040-dee2-421b-8e4f-1d316658501e.md")
$ws2.Hyperlinks.Add($ws2.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/276dd498861d647be956936a380962ede6cdd123/e2e/8fc7b040-dee2-421b-8e4f-1d316658501e.md", "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a3ee5f12757327b9da8f46c4f65da54df7a558c8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8fc7b040-dee2-421b-8e4f-1d316658501e.202cae47a453ea5feb7fbbec71dc0f8dbb40a093.zh-cn.xlf", "", "", "8fc7b040-dee2-421b-8e4f-1d316658501e.202cae47a453ea5feb7fbbec71dc0f8dbb40a093.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/59f5c454998e73de92b74483728bf5a6b20d5ec5/e2e/8fc7b040-dee2-421b-8e4f-1d316658501e.md", "", "", "8fc7b040-dee2-421b-8e4f-1d316658501e.md")
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2dd81b7c98e5f1723ae2e4ac576ceec3b8e7ec95/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8fc7b040-dee2-421b-8e4f-1d316658501e.202cae47a453ea5feb7fbbec71dc0f8dbb40a093.zh-cn.xlf", "", "", "8fc7b040-dee2-421b-8e4f-1d316658501e.202cae47a453ea5feb7fbbec71dc0f8dbb40a093.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("E2").Value = "2016-03-19 12:36:45"

# Row 3 (the ba5b2a27... file) no longer belongs in the report.
$ws3.Rows.Item(3).Delete()

# Rebuild the hyperlink collection so only the row-2 links survive.
$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/276dd498861d647be956936a380962ede6cdd123/e2e/8fc7b040-dee2-421b-8e4f-1d316658501e.md", "", "", "8fc7b040-dee2-421b-8e4f-1d316658501e.md")
$ws3.Hyperlinks.Add($ws3.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/276dd498861d647be956936a380962ede6cdd123/e2e/8fc7b040-dee2-421b-8e4f-1d316658501e.md", "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2f90f271fa9e0e4ff2c77f6ced47d2ca515da077/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8fc7b040-dee2-421b-8e4f-1d316658501e.202cae47a453ea5feb7fbbec71dc0f8dbb40a093.de-de.xlf", "", "", "8fc7b040-dee2-421b-8e4f-1d316658501e.202cae47a453ea5feb7fbbec71dc0f8dbb40a093.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/4b75385ddccb7c6c23d209e2cc43c6db95de75f0/e2e/8fc7b040-dee2-421b-8e4f-1d316658501e.md", "", "", "8fc7b040-dee2-421b-8e4f-1d316658501e.md")
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1ddf0acf5f00f2f9500f1d9c957fbc57b15470ac/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8fc7b040-dee2-421b-8e4f-1d316658501e.202cae47a453ea5feb7fbbec71dc0f8dbb40a093.de-de.xlf", "", "", "8fc7b040-dee2-421b-8e4f-1d316658501e.202cae47a453ea5feb7fbbec71dc0f8dbb40a093.de-de.xlf")

$wb.Save()
